$wb = $excel.ActiveWorkbook

# --- Sheet1: Train Results ---
$ws1 = $wb.Worksheets.Item("Train Results")

$data1 = @(
  @(0,40,4,0,28,4,4,20,2.31,2.555257558822632),
  @(4,0,8,24,4,16,44,0,2.98,2.986669063568115),
  @(4,16,0,20,4,12,44,0,3.22,3.12506890296936),
  @(0,24,4,0,8,4,36,24,2.45,2.452155828475952),
  @(4,0,8,4,24,12,44,4,3.17,2.967180967330933),
  @(0,16,8,4,16,20,36,0,3.07,3.021918296813965),
  @(4,20,4,4,16,0,52,0,3.35,3.11779522895813),
  @(4,0,12,4,4,16,52,8,2.81,2.800749063491821),
  @(4,12,8,0,8,8,56.00000000000001,4,2.88,2.876400947570801),
  @(4,4,4,8,28,16,32,4,3.04,3.049027919769287),
  @(4,0,8,4,24,12,44,4,3.09,2.967180967330933),
  @(4,8,0,12,20,8,48,0,3.11,3.136422872543335),
  @(0,20,8,4,0,20,44,4,2.82,2.991344213485718),
  @(0,20,8,4,0,20,44,4,2.74,2.991344213485718),
  @(4,12,0,0,16,8,52,8,2.62,2.885838508605957),
  @(4,12,8,4,16,8,48,0,3.19,3.045923471450806),
  @(0,20,4,0,4,4,48,20,2.52,2.446153879165649),
  @(4,8,4,4,16,12,48,4,2.94,2.999975442886353),
  @(4,12,0,12,20,8,44,0,3.29,3.156802892684937),
  @(4,12,0,12,20,8,44,0,3.25,3.156802892684937),
  @(0,12,8,4,4,20,36,4,3.01,2.997825145721436),
  @(4,12,4,4,20,16,32,8,2.96,2.948690414428711),
  @(4,8,0,0,4,8,52,24,2.44,2.409286975860596),
  @(4,0,4,0,20,8,52,12,2.92,2.741132736206055),
  @(4,8,12,4,4,24,39.99999999999999,4,2.9,2.912856817245483),
  @(4,0,12,16,4,12,52,0,2.96,2.912455081939697),
  @(0,16,8,4,16,20,36,0,3.02,3.021918296813965),
  @(4,12,8,0,8,8,56.00000000000001,4,2.79,2.876400947570801),
  @(4,0,4,4,16,12,56.00000000000001,4,2.85,2.949710845947266),
  @(0,4,8,4,16,20,44,4,2.94,2.960021018981934),
  @(0,12,4,0,4,16,39.99999999999999,24,2.51,2.649425506591797),
  @(0,0,8,4,16,12,52,8,2.99,2.896315336227417),
  @(0,16,0,4,20,20,28,12,3.38,3.021301984786987),
  @(4,16,8,0,12,0,48,12,2.56,2.592084407806396),
  @(0,0,8,4,16,12,52,8,2.82,2.896315336227417),
  @(4,12,4,4,12,16,36,12,2.86,2.852360486984253),
  @(0,24,8,0,12,16,32,8,2.93,2.930969715118408),
  @(4,32,8,4,4,20,24,4,2.84,3.019494533538818),
  @(0,24,8,0,12,16,32,8,2.94,2.930969715118408),
  @(4,40,0,0,12,4,36,4,3.16,3.060198545455933),
  @(0,16,0,4,20,20,28,12,2.72,3.021301984786987),
  @(4,20,4,4,16,0,52,0,3.21,3.11779522895813),
  @(4,12,0,0,16,8,52,8,2.56,2.885838508605957)
)

for ($i = 0; $i -lt $data1.Length; $i++) {
    $row = $data1[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws1.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}

# --- Sheet2: Test Results ---
$ws2 = $wb.Worksheets.Item("Test Results")

$data2 = @(
  @(0,8,0,4,4,28,36,20,2.56,2.898811101913452),
  @(4,16,8,0,12,0,48,12,2.67,2.592084407806396),
  @(4,16,0,4,12,8,52,4,2.96,3.028021097183228),
  @(4,8,4,4,16,12,48,4,2.92,2.999975442886353),
  @(4,0,4,0,4,4,60,24,2.45,2.380550861358643),
  @(4,0,4,0,20,8,52,12,2.64,2.741132736206055),
  @(4,12,4,0,28,0,39.99999999999999,12,2.94,2.790494441986084),
  @(4,4,12,0,0,4,52,24,2.22,2.39924693107605),
  @(4,12,8,4,16,8,48,0,3.15,3.045923233032227)
)

for ($i = 0; $i -lt $data2.Length; $i++) {
    $row = $data2[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws2.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}
